$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2005012531328321
$ws.Range("C2").Value = 0.5488721804511278
$ws.Range("J2").Value = 0.005012531328320802
$ws.Range("P2").Value = 0.1453634085213033
$ws.Range("S2").Value = 0.100250626566416
$ws.Range("C3").Value = 0.004608294930875576
$ws.Range("J3").Value = 0.02764976958525346
$ws.Range("P3").Value = 0.783410138248848
$ws.Range("S3").Value = 0.184331797235023
$ws.Range("J4").Value = 0.02
$ws.Range("P4").Value = 0.68
$ws.Range("S4").Value = 0.3
$ws.Range("B6").Value = 0.06334841628959276
$ws.Range("D6").Value = 0.01357466063348416
$ws.Range("E6").Value = 0.004524886877828055
$ws.Range("F6").Value = 0.04524886877828054
$ws.Range("J6").Value = 0.2624434389140272
$ws.Range("O6").Value = 0.009049773755656109
$ws.Range("Q6").Value = 0.1266968325791855
$ws.Range("R6").Value = 0.1176470588235294
$ws.Range("S6").Value = 0.3574660633484163
$ws.Range("B7").Value = 0.125
$ws.Range("D7").Value = 0.015
$ws.Range("F7").Value = 0.02
$ws.Range("J7").Value = 0.13
$ws.Range("O7").Value = 0.015
$ws.Range("Q7").Value = 0.14
$ws.Range("R7").Value = 0.075
$ws.Range("S7").Value = 0.48
$ws.Range("B8").Value = 0.1393258426966292
$ws.Range("D8").Value = 0.01797752808988764
$ws.Range("E8").Value = 0.006741573033707865
$ws.Range("F8").Value = 0.05168539325842696
$ws.Range("J8").Value = 0.1056179775280899
$ws.Range("O8").Value = 0.02022471910112359
$ws.Range("Q8").Value = 0.09438202247191012
$ws.Range("R8").Value = 0.09662921348314607
$ws.Range("S8").Value = 0.4674157303370787
$ws.Range("B9").Value = 0.1216931216931217
$ws.Range("D9").Value = 0.01058201058201058
$ws.Range("F9").Value = 0.04232804232804233
$ws.Range("J9").Value = 0.1111111111111111
$ws.Range("O9").Value = 0.01587301587301587
$ws.Range("Q9").Value = 0.1428571428571428
$ws.Range("R9").Value = 0.08994708994708994
$ws.Range("S9").Value = 0.4656084656084656
$ws.Range("B10").Value = 0.1308864265927978
$ws.Range("D10").Value = 0.02562326869806094
$ws.Range("F10").Value = 0.07479224376731301
$ws.Range("J10").Value = 0.1177285318559557
$ws.Range("O10").Value = 0.0131578947368421
$ws.Range("Q10").Value = 0.1745152354570637
$ws.Range("R10").Value = 0.07548476454293629
$ws.Range("S10").Value = 0.3878116343490305
$ws.Range("G11").Value = 0.1335504885993485
$ws.Range("J11").Value = 0.07817589576547231
$ws.Range("K11").Value = 0.1889250814332248
$ws.Range("L11").Value = 0.5830618892508144
$ws.Range("S11").Value = 0.01628664495114007
$ws.Range("G12").Value = 0.6864864864864865
$ws.Range("J12").Value = 0.2324324324324324
$ws.Range("L12").Value = 0.02162162162162162
$ws.Range("S12").Value = 0.05945945945945946
$ws.Range("G13").Value = 0.7826086956521739
$ws.Range("J13").Value = 0.1956521739130435
$ws.Range("S13").Value = 0.02173913043478261
$ws.Range("F15").Value = 0.02392344497607655
$ws.Range("H15").Value = 0.1004784688995215
$ws.Range("I15").Value = 0.09569377990430622
$ws.Range("J15").Value = 0.354066985645933
$ws.Range("K15").Value = 0.07177033492822966
$ws.Range("O15").Value = 0.05263157894736842
$ws.Range("S15").Value = 0.3014354066985646
$ws.Range("F16").Value = 0.01945525291828794
$ws.Range("H16").Value = 0.1595330739299611
$ws.Range("I16").Value = 0.08560311284046693
$ws.Range("J16").Value = 0.4396887159533074
$ws.Range("K16").Value = 0.07782101167315175
$ws.Range("M16").Value = 0.01556420233463035
$ws.Range("O16").Value = 0.03501945525291829
$ws.Range("S16").Value = 0.1673151750972763
$ws.Range("F17").Value = 0.01851851851851852
$ws.Range("H17").Value = 0.1455026455026455
$ws.Range("I17").Value = 0.06084656084656084
$ws.Range("J17").Value = 0.5132275132275133
$ws.Range("K17").Value = 0.08465608465608465
$ws.Range("M17").Value = 0.01058201058201058
$ws.Range("O17").Value = 0.04497354497354497
$ws.Range("S17").Value = 0.1216931216931217
$ws.Range("F18").Value = 0.01435406698564593
$ws.Range("H18").Value = 0.1913875598086124
$ws.Range("I18").Value = 0.09569377990430622
$ws.Range("J18").Value = 0.3684210526315789
$ws.Range("K18").Value = 0.08133971291866028
$ws.Range("M18").Value = 0.02870813397129187
$ws.Range("O18").Value = 0.03827751196172249
$ws.Range("S18").Value = 0.1818181818181818
$ws.Range("F19").Value = 0.01566579634464752
$ws.Range("H19").Value = 0.1912532637075718
$ws.Range("I19").Value = 0.06919060052219321
$ws.Range("J19").Value = 0.3831592689295039
$ws.Range("K19").Value = 0.1011749347258486
$ws.Range("M19").Value = 0.02349869451697128
$ws.Range("N19").Value = 0.001305483028720627
$ws.Range("O19").Value = 0.06527415143603134
$ws.Range("S19").Value = 0.1494778067885117
